$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(3).Insert()

$ws.Cells.Item(1, 3).Value = "statut_name"

$ws.Cells.Item(2, 3).Value = "résultat et / ou publication posté"
$ws.Cells.Item(3, 3).Value = "résultat et / ou publication posté dans les 12 mois"
$ws.Cells.Item(4, 3).Value = "résultat et / ou publication posté dans les 36 mois"
$ws.Cells.Item(5, 3).Value = "pas de résultat ni de publication"
$ws.Cells.Item(6, 3).Value = "pas de résultat ni de publication"
$ws.Cells.Item(7, 3).Value = "pas de résultat ni de publication"
